$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "D", "E", "H", "J")
$replacements = @(
    @("D64", "D69"),
    @("D51", "D55"),
    @("D80", "D86"),
    @("S30", "S31")
)

$lastRow = $ws.UsedRange.Rows.Count

foreach ($col in $cols) {
    $rng = $ws.Range($col + "1:" + $col + $lastRow)
    foreach ($pair in $replacements) {
        $rng.Replace($pair[0], $pair[1]) | Out-Null
    }
}
